# Update weekly Fruta/Hortaliza price records: the rows' date/volume/price/origin
# data (columns D, M, N, O, P, R, S) have been reshuffled among rows 2-12 and 14-19
# (row 13 is unchanged). Row 1 (headers) and columns A, B, C, E, F, G, H, I, J, K,
# L, Q, T are identical for every data row, so only the moved columns are touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44589
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 6000
$ws.Range("O2").Value = 6000
$ws.Range("P2").Value = 6000
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 3000

$ws.Range("D3").Value = 44214
$ws.Range("M3").Value = 48
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 6000
$ws.Range("P3").Value = 6000
$ws.Range("R3").Value = "Provincia de Linares"
$ws.Range("S3").Value = 3000

$ws.Range("D4").Value = 44592
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = 8000
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 8000
$ws.Range("R4").Value = "Provincia de Linares"
$ws.Range("S4").Value = 4000

$ws.Range("D5").Value = 44606
$ws.Range("M5").Value = 45
$ws.Range("N5").Value = 7000
$ws.Range("O5").Value = 7000
$ws.Range("P5").Value = 7000
$ws.Range("R5").Value = "Provincia de Linares"
$ws.Range("S5").Value = 3500

$ws.Range("D6").Value = 44974
$ws.Range("M6").Value = 130
$ws.Range("N6").Value = 7000
$ws.Range("O6").Value = 7500
$ws.Range("P6").Value = 7269
$ws.Range("R6").Value = "Provincia de Curicó"
$ws.Range("S6").Value = 3634

$ws.Range("D7").Value = 44209
$ws.Range("M7").Value = 58
$ws.Range("N7").Value = 6000
$ws.Range("O7").Value = 6000
$ws.Range("P7").Value = 6000
$ws.Range("R7").Value = "Provincia de Curicó"
$ws.Range("S7").Value = 3000

$ws.Range("D8").Value = 44614
$ws.Range("M8").Value = 45
$ws.Range("N8").Value = 6000
$ws.Range("O8").Value = 6000
$ws.Range("P8").Value = 6000
$ws.Range("R8").Value = "Provincia de Linares"
$ws.Range("S8").Value = 3000

$ws.Range("D9").Value = 44627
$ws.Range("M9").Value = 45
$ws.Range("N9").Value = 6000
$ws.Range("O9").Value = 6000
$ws.Range("P9").Value = 6000
$ws.Range("R9").Value = "Provincia de Linares"
$ws.Range("S9").Value = 3000

$ws.Range("D10").Value = 44628
$ws.Range("M10").Value = 40
$ws.Range("N10").Value = 6000
$ws.Range("O10").Value = 6000
$ws.Range("P10").Value = 6000
$ws.Range("R10").Value = "Provincia de Linares"
$ws.Range("S10").Value = 3000

$ws.Range("D11").Value = 44960
$ws.Range("M11").Value = 40
$ws.Range("N11").Value = 7000
$ws.Range("O11").Value = 7000
$ws.Range("P11").Value = 7000
$ws.Range("R11").Value = "Provincia de Curicó"
$ws.Range("S11").Value = 3500

$ws.Range("D12").Value = 44211
$ws.Range("M12").Value = 45
$ws.Range("N12").Value = 6000
$ws.Range("O12").Value = 6000
$ws.Range("P12").Value = 6000
$ws.Range("R12").Value = "Provincia de Curicó"
$ws.Range("S12").Value = 3000

$ws.Range("D14").Value = 44959
$ws.Range("M14").Value = 40
$ws.Range("N14").Value = 7000
$ws.Range("O14").Value = 7000
$ws.Range("P14").Value = 7000
$ws.Range("R14").Value = "Provincia de Curicó"
$ws.Range("S14").Value = 3500

$ws.Range("D15").Value = 44586
$ws.Range("M15").Value = 80
$ws.Range("N15").Value = 7000
$ws.Range("O15").Value = 7000
$ws.Range("P15").Value = 7000
$ws.Range("R15").Value = "Provincia de Curicó"
$ws.Range("S15").Value = 3500

$ws.Range("D16").Value = 44587
$ws.Range("M16").Value = 165
$ws.Range("N16").Value = 6500
$ws.Range("O16").Value = 7000
$ws.Range("P16").Value = 6742
$ws.Range("R16").Value = "Provincia de Linares"
$ws.Range("S16").Value = 3371

$ws.Range("D17").Value = 45001
$ws.Range("M17").Value = 66
$ws.Range("N17").Value = 7500
$ws.Range("O17").Value = 8000
$ws.Range("P17").Value = 7773
$ws.Range("R17").Value = "Provincia de Curicó"
$ws.Range("S17").Value = 3886

$ws.Range("D18").Value = 44588
$ws.Range("M18").Value = 160
$ws.Range("N18").Value = 6500
$ws.Range("O18").Value = 7000
$ws.Range("P18").Value = 6750
$ws.Range("R18").Value = "Provincia de Curicó"
$ws.Range("S18").Value = 3375

$ws.Range("D19").Value = 44585
$ws.Range("M19").Value = 160
$ws.Range("N19").Value = 6500
$ws.Range("O19").Value = 7000
$ws.Range("P19").Value = 6750
$ws.Range("R19").Value = "Provincia de Curicó"
$ws.Range("S19").Value = 3375
